$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextCell "D2" "62.251.09"
Set-TextCell "E2" "  +0.51%  "
Set-TextCell "D3" "3.427.64"
Set-TextCell "E3" "  +0.26%  "
Set-TextCell "D4" "1.00"
Set-TextCell "E4" "  -0.34%  "
Set-TextCell "D5" "414.05"
Set-TextCell "E5" "  +1.07%  "
Set-TextCell "D6" "129.19"
Set-TextCell "E6" "  +0.68%  "
Set-TextCell "D7" "0.623"
Set-TextCell "E7" "  -1.61%  "
Set-TextCell "D8" "1.00"
Set-TextCell "E8" "  +0.06%  "
Set-TextCell "D9" "0.725"
Set-TextCell "E9" "  -0.66%  "
Set-TextCell "E10" "  +0.81%  "
Set-TextCell "D11" "42.76"
Set-TextCell "E11" "  +0.53%  "
Set-TextCell "D12" "9.30"
Set-TextCell "E12" "  +3.05%  "
Set-TextCell "B13" "ShibaInu"
Set-TextCell "C13" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D13" "0.0000216"
Set-TextCell "E13" "  +6.11%  "
Set-TextCell "B14" "WrappedliquidstakedEther2.0"
Set-TextCell "C14" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell "D14" "3.972.14"
Set-TextCell "E14" "  +0.35%  "
Set-TextCell "D16" "20.50"
Set-TextCell "E16" "  -3.25%  "
Set-TextCell "D17" "3.452.44"
Set-TextCell "E17" "  +1.52%  "
Set-TextCell "D18" "12.59"
Set-TextCell "E18" "  +2.55%  "
Set-TextCell "D19" "1.07"
Set-TextCell "E19" "  -0.54%  "
Set-TextCell "D20" "62.283.14"
Set-TextCell "E20" "  +0.68%  "
Set-TextCell "D21" "465.92"
Set-TextCell "E21" "  +4.52%  "
Set-TextCell "D22" "90.89"
Set-TextCell "E22" "  -1.39%  "
Set-TextCell "E23" "  +3.60%  "
Set-TextCell "D24" "13.54"
Set-TextCell "E24" "  +5.04%  "
Set-TextCell "D25" "10.54"
Set-TextCell "E25" "  +19.68%  "
Set-TextCell "D26" "3.31"
Set-TextCell "E26" "  +2.34%  "
Set-TextCell "D27" "33.14"
Set-TextCell "E27" "  -0.13%  "
Set-TextCell "E28" "  -0.30%  "
Set-TextCell "D29" "7.68"
Set-TextCell "E29" "  +1.73%  "
Set-TextCell "D30" "11.96"
Set-TextCell "E30" "  +0.11%  "
Set-TextCell "E31" "  -3.35%  "
Set-TextCell "E32" "  -1.41%  "
Set-TextCell "E33" "  -1.60%  "
Set-TextCell "D34" "40.68"
Set-TextCell "E34" "  -4.67%  "
Set-TextCell "E35" "  +0.06%  "
Set-TextCell "D36" "58.74"
Set-TextCell "E36" "  +9.98%  "
Set-TextCell "D37" "0.0487"
Set-TextCell "E37" "  -1.57%  "
Set-TextCell "D38" "0.999"
Set-TextCell "E38" "  +0.05%  "
Set-TextCell "E39" "  +5.11%  "
Set-TextCell "D40" "0.326"
Set-TextCell "E40" "  +4.09%  "
Set-TextCell "B41" "Stellar"
Set-TextCell "C41" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D41" "0.134"
Set-TextCell "E41" "  +0.27%  "
Set-TextCell "B42" "LidoDAOToken"
Set-TextCell "C42" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D42" "3.34"
Set-TextCell "E42" "  -1.01%  "
Set-TextCell "D43" "145.39"
Set-TextCell "E43" "  +3.40%  "
Set-TextCell "E44" "  +10.18%  "
Set-TextCell "E45" "  +5.51%  "
Set-TextCell "D46" "4.31"
Set-TextCell "E46" "  +2.16%  "
Set-TextCell "E47" "  +19.89%  "
Set-TextCell "D48" "16.45"
Set-TextCell "E48" "  -0.19%  "
Set-TextCell "D49" "22.29"
Set-TextCell "E49" "  -0.62%  "
Set-TextCell "D50" "0.0₃0518"
Set-TextCell "E50" "  +26.56%  "
Set-TextCell "D51" "110.46"
Set-TextCell "E51" "  +5.97%  "
